$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -13.312
$ws.Range("B3").Value = 5.914
$ws.Range("C3").Value = -12.244
$ws.Range("B4").Value = 6.752999999999998
$ws.Range("E8").Value = 16.679
$ws.Range("C9").Value = -10.988
$ws.Range("A11").Value = -21.518
$ws.Range("E11").Value = 16.712
$ws.Range("A12").Value = -21.654
$ws.Range("B14").Value = 5.865
$ws.Range("E14").Value = 16.989
$ws.Range("A15").Value = -21.672
$ws.Range("C15").Value = -13.258
$ws.Range("E15").Value = 16.175
$ws.Range("E17").Value = 16.494
$ws.Range("C19").Value = -12.869
$ws.Range("C20").Value = -12.219
$ws.Range("C25").Value = -12.556
$ws.Range("B26").Value = 6.144
$ws.Range("E26").Value = 16.388
$ws.Range("A27").Value = -21.313
$ws.Range("C27").Value = -13.685
$ws.Range("A28").Value = -21.463
$ws.Range("C28").Value = -13.334
$ws.Range("C30").Value = -13.481
$ws.Range("A31").Value = -21.106
$ws.Range("B31").Value = 6.157000000000001
$ws.Range("A32").Value = -21.586
$ws.Range("C32").Value = -12.717
$ws.Range("B35").Value = 6.726000000000001
$ws.Range("A36").Value = -21.071
$ws.Range("E36").Value = 16.868
$ws.Range("B37").Value = 7.113
$ws.Range("A38").Value = -20.108
$ws.Range("B39").Value = 7.586
$ws.Range("B40").Value = 8.571
$ws.Range("E42").Value = 16.419
$ws.Range("C44").Value = -12.611
$ws.Range("B45").Value = 5.612
$ws.Range("A46").Value = -21.72
$ws.Range("C47").Value = -12.321
$ws.Range("B52").Value = 4.798
$ws.Range("A54").Value = -21.28100000000001
$ws.Range("A55").Value = -21.825
$ws.Range("A56").Value = -21.538
$ws.Range("B57").Value = 5.827999999999999
$ws.Range("C58").Value = -12.916
$ws.Range("C62").Value = -13.329
$ws.Range("E64").Value = 17.25
$ws.Range("A67").Value = -21.6
$ws.Range("E68").Value = 17.078
$ws.Range("A69").Value = -21.651
$ws.Range("A72").Value = -21.567
$ws.Range("A73").Value = -20.552
$ws.Range("C77").Value = -13.51
$ws.Range("C78").Value = -13.679
$ws.Range("E79").Value = 17.312
$ws.Range("B81").Value = 6.574
$ws.Range("A83").Value = -20.489
$ws.Range("B83").Value = 6.728
$ws.Range("C84").Value = -13.476
$ws.Range("A86").Value = -22.076
$ws.Range("C89").Value = -11.302
$ws.Range("E89").Value = 16.975
$ws.Range("A91").Value = -21.768
$ws.Range("C91").Value = -11.534
$ws.Range("C92").Value = -11.465
$ws.Range("A93").Value = -21.593
$ws.Range("C96").Value = -13.343
$ws.Range("A99").Value = -21.067
$ws.Range("B100").Value = 5.789
$ws.Range("B102").Value = 6.214
$ws.Range("C102").Value = -12.686
